# New weekly price record for Papa (Patagonia, "1a (guarda)") reported at
# Terminal Hortofrutícola Agro Chillán is inserted as row 175. Inserting the
# row (rather than just appending) pushes every existing record from the old
# row 175 down through the old row 201 one row further down (new rows 176-202),
# matching the target workbook exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 175; this shifts rows 175-201 down to 176-202 and
# extends the sheet dimension from A1:R201 to A1:R202 automatically.
$ws.Rows.Item(175).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(175, 1).Value = 7
$ws.Cells.Item(175, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(175, 3).Value = "Ñuble"
$ws.Cells.Item(175, 4).Value = 44474
$ws.Cells.Item(175, 5).Value = 16
$ws.Cells.Item(175, 6).Value = 100114001
$ws.Cells.Item(175, 7).Value = "Papa"
$ws.Cells.Item(175, 8).Value = "Patagonia"
$ws.Cells.Item(175, 9).Value = "1a (guarda)"
$ws.Cells.Item(175, 10).Value = 120
$ws.Cells.Item(175, 11).Value = 6500
$ws.Cells.Item(175, 12).Value = 7000
$ws.Cells.Item(175, 13).Value = 6750
$ws.Cells.Item(175, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(175, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(175, 16).Value = 270
$ws.Cells.Item(175, 17).Value = 25
$ws.Cells.Item(175, 18).Value = "Hortaliza"
